$wb = $excel.ActiveWorkbook

# Sheet ALC
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H43").Value = 7749.75
$ws.Range("J43").Value = 3999
$ws.Range("L43").Value = 3999
$ws.Range("N43").Value = -4137
$ws.Range("H62").Value = 4000
$ws.Range("I62").Value = 0
$ws.Range("J62").Value = 4000
$ws.Range("K62").Value = 0
$ws.Range("L62").Value = 4000
$ws.Range("M62").ClearContents()
$ws.Range("N62").Value = -5248
$ws.Range("H65").Value = 4000
$ws.Range("I65").Value = 0
$ws.Range("J65").Value = 4000
$ws.Range("K65").Value = 0
$ws.Range("L65").Value = 20000
$ws.Range("M65").ClearContents()
$ws.Range("N65").Value = -26240
$ws.Range("H69").Value = 16998.117
$ws.Range("J69").Value = 18499.857
$ws.Range("L69").Value = 55499.571
$ws.Range("N69").Value = -57247.571
$ws.Range("H72").Value = 16998.117
$ws.Range("J72").Value = 18499.857
$ws.Range("L72").Value = 166498.713
$ws.Range("N72").Value = -175234.713
$ws.Range("H74").Value = 2595.6
$ws.Range("I74").Value = 2595.6
$ws.Range("J74").Value = 0
$ws.Range("K74").Value = 2595.6
$ws.Range("L74").Value = 0
$ws.Range("M74").Value = -1659.6
$ws.Range("N74").ClearContents()
$ws.Range("H77").Value = 2595.6
$ws.Range("I77").Value = 2595.6
$ws.Range("J77").Value = 0
$ws.Range("K77").Value = 12978
$ws.Range("L77").Value = 0
$ws.Range("M77").Value = -8298
$ws.Range("N77").ClearContents()
$ws.Range("H88").Value = 26367510
$ws.Range("I88").Value = 166668130
$ws.Range("J88").Value = 2984075.2
$ws.Range("K88").Value = 166668130
$ws.Range("L88").Value = 2984075.2
$ws.Range("M88").Value = -166667724
$ws.Range("N88").Value = -2984887.2
$ws.Range("H91").Value = 26367510
$ws.Range("I91").Value = 166668130
$ws.Range("J91").Value = 2984075.2
$ws.Range("K91").Value = 166668130
$ws.Range("L91").Value = 2984075.2
$ws.Range("M91").Value = -166666726
$ws.Range("N91").Value = -2986883.2
$ws.Range("H103").Value = 1312.5714
$ws.Range("J103").Value = 1147
$ws.Range("L103").Value = 3441
$ws.Range("N103").Value = -4613
$ws.Range("H125").Value = 10928204
$ws.Range("I125").Value = 5651549
$ws.Range("J125").Value = 12511200
$ws.Range("K125").Value = 50863941
$ws.Range("L125").Value = 112600800
$ws.Range("M125").Value = -50861481
$ws.Range("N125").Value = -112605720
$ws.Range("H132").Value = 4491.4546
$ws.Range("I132").Value = 4491.4546
$ws.Range("K132").Value = 13474.3638
$ws.Range("M132").Value = -10944.3638

# Sheet ARM
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 5407.241
$ws.Range("I32").Value = 2547.7917
$ws.Range("K32").Value = 2547.7917
$ws.Range("M32").Value = -2260.7917
$ws.Range("H61").Value = 52633344
$ws.Range("I61").Value = 58825296
$ws.Range("J61").Value = 1749
$ws.Range("K61").Value = 58825296
$ws.Range("L61").Value = 1749
$ws.Range("M61").Value = -58825084
$ws.Range("N61").Value = -2173
$ws.Range("H63").Value = 5666.6665
$ws.Range("J63").Value = 2000
$ws.Range("L63").Value = 2000
$ws.Range("N63").Value = -3372
$ws.Range("H66").Value = 5666.6665
$ws.Range("J66").Value = 2000
$ws.Range("L66").Value = 10000
$ws.Range("N66").Value = -16864
$ws.Range("H132").Value = 6252944
$ws.Range("I132").Value = 6669473.5
$ws.Range("K132").Value = 20008420.5
$ws.Range("M132").Value = -20005890.5
$ws.Range("H136").Value = 52633344
$ws.Range("I136").Value = 58825296
$ws.Range("J136").Value = 1749
$ws.Range("K136").Value = 176475888
$ws.Range("L136").Value = 5247
$ws.Range("M136").Value = -176473338
$ws.Range("N136").Value = -10347

# Sheet BSM
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H20").Value = 2969.1667
$ws.Range("I20").Value = 3059.7
$ws.Range("J20").Value = 2856
$ws.Range("K20").Value = 3059.7
$ws.Range("L20").Value = 2856
$ws.Range("M20").Value = -2812.7
$ws.Range("N20").Value = -3350
$ws.Range("H76").Value = 54999
$ws.Range("J76").Value = 54999
$ws.Range("L76").Value = 54999
$ws.Range("N76").Value = -55629
$ws.Range("H79").Value = 54999
$ws.Range("J79").Value = 54999
$ws.Range("L79").Value = 54999
$ws.Range("N79").Value = -57183
$ws.Range("H107").Value = 49920.24
$ws.Range("I107").Value = 2501.4443
$ws.Range("J107").Value = 334433
$ws.Range("K107").Value = 2501.4443
$ws.Range("L107").Value = 334433
$ws.Range("M107").Value = -581.4443000000001
$ws.Range("N107").Value = -338273

# Sheet CRP
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 8821.75
$ws.Range("I31").Value = 6139.4443
$ws.Range("K31").Value = 6139.4443
$ws.Range("M31").Value = -5844.4443
$ws.Range("H34").Value = 8821.75
$ws.Range("I34").Value = 6139.4443
$ws.Range("K34").Value = 6139.4443
$ws.Range("M34").Value = -5937.4443
$ws.Range("H58").Value = 23816038
$ws.Range("I58").Value = 38470684
$ws.Range("J58").Value = 2239.75
$ws.Range("K58").Value = 38470684
$ws.Range("L58").Value = 2239.75
$ws.Range("M58").Value = -38470481
$ws.Range("N58").Value = -2645.75
$ws.Range("H132").Value = 111116104
$ws.Range("I132").Value = 200004590
$ws.Range("K132").Value = 600013770
$ws.Range("M132").Value = -600011240
$ws.Range("H136").Value = 23816038
$ws.Range("I136").Value = 38470684
$ws.Range("J136").Value = 2239.75
$ws.Range("K136").Value = 115412052
$ws.Range("L136").Value = 6719.25
$ws.Range("M136").Value = -115409502
$ws.Range("N136").Value = -11819.25

# Sheet CUL
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H32").Value = 33556244
$ws.Range("J32").Value = 334066
$ws.Range("L32").Value = 1002198
$ws.Range("N32").Value = -1002764
$ws.Range("H122").Value = 1058.75
$ws.Range("I122").Value = 982.8182
$ws.Range("K122").Value = 8845.363800000001
$ws.Range("M122").Value = -6395.363800000001
$ws.Range("H125").Value = 9999
$ws.Range("J125").Value = 9999
$ws.Range("L125").Value = 29997
$ws.Range("N125").Value = -39837

# Sheet GSM
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H2").Value = 106.333336
$ws.Range("I2").Value = 36.77778
$ws.Range("J2").Value = 210.66667
$ws.Range("K2").Value = 36.77778
$ws.Range("L2").Value = 210.66667
$ws.Range("M2").Value = 76.22221999999999
$ws.Range("N2").Value = -436.66667
$ws.Range("H80").Value = 1500
$ws.Range("I80").Value = 2000
$ws.Range("J80").Value = 1000
$ws.Range("K80").Value = 2000
$ws.Range("L80").Value = 1000
$ws.Range("M80").Value = -1002
$ws.Range("N80").Value = -2996
$ws.Range("H83").Value = 1500
$ws.Range("I83").Value = 2000
$ws.Range("J83").Value = 1000
$ws.Range("K83").Value = 10000
$ws.Range("L83").Value = 5000
$ws.Range("M83").Value = -5008
$ws.Range("N83").Value = -14984
$ws.Range("H113").Value = 42568.8
$ws.Range("I113").Value = 58248.668
$ws.Range("J113").Value = 2249.1428
$ws.Range("K113").Value = 58248.668
$ws.Range("L113").Value = 2249.1428
$ws.Range("M113").Value = -56078.668
$ws.Range("N113").Value = -6589.1428
$ws.Range("H122").Value = 58605.953
$ws.Range("I122").Value = 64806.895
$ws.Range("K122").Value = 194420.685
$ws.Range("M122").Value = -191970.685
$ws.Range("H132").Value = 5437403
$ws.Range("I132").Value = 5684422
$ws.Range("J132").Value = 2982
$ws.Range("K132").Value = 17053266
$ws.Range("L132").Value = 8946
$ws.Range("M132").Value = -17050736
$ws.Range("N132").Value = -14006

# Sheet LTW
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H16").Value = 1055.9375
$ws.Range("I16").Value = 793
$ws.Range("K16").Value = 793
$ws.Range("M16").Value = -623
$ws.Range("H22").Value = 3053.5454
$ws.Range("I22").Value = 3849.8333
$ws.Range("J22").Value = 2098
$ws.Range("K22").Value = 3849.8333
$ws.Range("L22").Value = 2098
$ws.Range("M22").Value = -3554.8333
$ws.Range("N22").Value = -2688
$ws.Range("H27").Value = 3053.5454
$ws.Range("I27").Value = 3849.8333
$ws.Range("J27").Value = 2098
$ws.Range("K27").Value = 3849.8333
$ws.Range("L27").Value = 2098
$ws.Range("M27").Value = -3742.8333
$ws.Range("N27").Value = -2312
$ws.Range("H82").Value = 1369.4
$ws.Range("I82").Value = 1369.4
$ws.Range("K82").Value = 1369.4
$ws.Range("M82").Value = -1008.4
$ws.Range("H85").Value = 1369.4
$ws.Range("I85").Value = 1369.4
$ws.Range("K85").Value = 1369.4
$ws.Range("M85").Value = -121.4000000000001

# Sheet WVR
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H126").Value = 2139.3157
$ws.Range("I126").Value = 2639.6667
$ws.Range("K126").Value = 7919.000100000001
$ws.Range("M126").Value = -5449.000100000001
